# Apply the Arabic translation edits described by the diff.
$d = $word.ActiveDocument

# 1. "{Onboarding}" -> "{Onboarding} " (title heading, add trailing space)
$d.Content.Find.Execute("{Onboarding}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{Onboarding} ", 2) | Out-Null

# 2. "{Lesson: Onboarding}" -> "{Lesson: Onboarding} " (table header, add trailing space)
$d.Content.Find.Execute("{Lesson: Onboarding}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{Lesson: Onboarding} ", 2) | Out-Null

# 3. "... MENU ..." -> "... قائمة ..." (keyword translated from English to Arabic)
$d.Content.Find.Execute("كتابة MENU والانتقال", $true, $false, $false, $false, $false,
                         $true, 1, $false, "كتابة قائمة والانتقال", 2) | Out-Null

# 4. "... HELP ..." -> "... مساعدة ..."
$d.Content.Find.Execute("كتابة HELP في أي وقت", $true, $false, $false, $false, $false,
                         $true, 1, $false, "كتابة مساعدة في أي وقت", 2) | Out-Null

# 5. "... PLAY ..." -> "... تشغيل ..."
$d.Content.Find.Execute("كتابة PLAY في أي وقت", $true, $false, $false, $false, $false,
                         $true, 1, $false, "كتابة تشغيل في أي وقت", 2) | Out-Null

# 6. "هادى ء" -> "هادي " — only the bare list-item occurrence (immediately followed
#    by a textWrapping line break), not the "هادى ء- للاسترخاء معًا، أو" phrase
#    elsewhere in the document which must stay untouched.
$needle = "هادى ء"
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
$targetStart = -1
$targetEnd = -1
while ($found) {
    $afterStart = $rng.End
    $afterEnd = [Math]::Min($afterStart + 1, $d.Content.End)
    $afterChar = $d.Range($afterStart, $afterEnd).Text
    if ($afterChar.Length -gt 0 -and [int][char]$afterChar[0] -eq 11) {
        $targetStart = $rng.Start
        $targetEnd = $rng.End
    }
    $nextStart = $rng.End
    $rng = $d.Range($nextStart, $d.Content.End)
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
}
if ($targetStart -ge 0) {
    $target = $d.Range($targetStart, $targetEnd)
    $target.Text = "هادي "
}
